$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 202
$ws.Range('D202').Value = 44706
$ws.Range('I202').Value = 'Especial'
$ws.Range('J202').Value = 220
$ws.Range('K202').Value = 14000
$ws.Range('L202').Value = 14000
$ws.Range('M202').Value = 14000
$ws.Range('P202').Value = 778

# Row 203
$ws.Range('D203').Value = 44706
$ws.Range('I203').Value = 'Primera'
$ws.Range('J203').Value = 250
$ws.Range('K203').Value = 12000
$ws.Range('L203').Value = 12000
$ws.Range('M203').Value = 12000
$ws.Range('P203').Value = 667

# Row 204
$ws.Range('D204').Value = 44706
$ws.Range('I204').Value = 'Segunda'
$ws.Range('J204').Value = 280
$ws.Range('K204').Value = 10000
$ws.Range('L204').Value = 10000
$ws.Range('M204').Value = 10000
$ws.Range('P204').Value = 556

# Row 205
$ws.Range('D205').Value = 44351
$ws.Range('I205').Value = 'Primera'
$ws.Range('J205').Value = 180
$ws.Range('K205').Value = 15000
$ws.Range('L205').Value = 15000
$ws.Range('M205').Value = 15000
$ws.Range('P205').Value = 833

# Row 206
$ws.Range('D206').Value = 44351
$ws.Range('I206').Value = 'Segunda'
$ws.Range('J206').Value = 260
$ws.Range('K206').Value = 13000
$ws.Range('L206').Value = 13000
$ws.Range('M206').Value = 13000
$ws.Range('P206').Value = 722

# Row 207
$ws.Range('D207').Value = 44351
$ws.Range('I207').Value = 'Tercera'
$ws.Range('J207').Value = 200
$ws.Range('K207').Value = 8500
$ws.Range('L207').Value = 8500
$ws.Range('M207').Value = 8500
$ws.Range('P207').Value = 472

# Row 208
$ws.Range('D208').Value = 44692
$ws.Range('J208').Value = 220
$ws.Range('N208').Value = '$/bandeja 18 kilos'

# Row 209
$ws.Range('D209').Value = 44692
$ws.Range('J209').Value = 250
$ws.Range('N209').Value = '$/bandeja 18 kilos'

# Row 210
$ws.Range('D210').Value = 44692
$ws.Range('J210').Value = 280
$ws.Range('N210').Value = '$/bandeja 18 kilos'

# Row 211
$ws.Range('D211').Value = 44315
$ws.Range('I211').Value = 'Especial'
$ws.Range('J211').Value = 350
$ws.Range('K211').Value = 14000
$ws.Range('L211').Value = 14000
$ws.Range('M211').Value = 14000
$ws.Range('N211').Value = '$/caja 18 kilos'
$ws.Range('P211').Value = 778

# Row 212
$ws.Range('D212').Value = 44315
$ws.Range('I212').Value = 'Primera'
$ws.Range('J212').Value = 580
$ws.Range('L212').Value = 12000
$ws.Range('M212').Value = 12000
$ws.Range('N212').Value = '$/caja 18 kilos'
$ws.Range('P212').Value = 667

# Row 213 (new)
$ws.Range('A213').Value = 12
$ws.Range('B213').Value = 'Mapocho Venta Directa de Santiago'
$ws.Range('C213').Value = 'Metropolitana'
$ws.Range('D213').Value = 44315
$ws.Range('D213').NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range('E213').Value = 13
$ws.Range('F213').Value = 100112043
$ws.Range('G213').Value = 'Pepino dulce'
$ws.Range('H213').Value = 'Cultivar IV Región'
$ws.Range('I213').Value = 'Segunda'
$ws.Range('J213').Value = 300
$ws.Range('K213').Value = 10000
$ws.Range('L213').Value = 10000
$ws.Range('M213').Value = 10000
$ws.Range('N213').Value = '$/caja 18 kilos'
$ws.Range('O213').Value = 'Provincia de Limarí'
$ws.Range('P213').Value = 556
$ws.Range('Q213').Value = 18
$ws.Range('R213').Value = 'Hortaliza'

# Row 214 (new)
$ws.Range('A214').Value = 12
$ws.Range('B214').Value = 'Mapocho Venta Directa de Santiago'
$ws.Range('C214').Value = 'Metropolitana'
$ws.Range('D214').Value = 44376
$ws.Range('D214').NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range('E214').Value = 13
$ws.Range('F214').Value = 100112043
$ws.Range('G214').Value = 'Pepino dulce'
$ws.Range('H214').Value = 'Cultivar IV Región'
$ws.Range('I214').Value = 'Primera'
$ws.Range('J214').Value = 550
$ws.Range('K214').Value = 15000
$ws.Range('L214').Value = 16000
$ws.Range('M214').Value = 15636
$ws.Range('N214').Value = '$/bandeja 18 kilos'
$ws.Range('O214').Value = 'Provincia de Limarí'
$ws.Range('P214').Value = 869
$ws.Range('Q214').Value = 18
$ws.Range('R214').Value = 'Hortaliza'

# Row 215 (new)
$ws.Range('A215').Value = 12
$ws.Range('B215').Value = 'Mapocho Venta Directa de Santiago'
$ws.Range('C215').Value = 'Metropolitana'
$ws.Range('D215').Value = 44376
$ws.Range('D215').NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range('E215').Value = 13
$ws.Range('F215').Value = 100112043
$ws.Range('G215').Value = 'Pepino dulce'
$ws.Range('H215').Value = 'Cultivar IV Región'
$ws.Range('I215').Value = 'Segunda'
$ws.Range('J215').Value = 390
$ws.Range('K215').Value = 12000
$ws.Range('L215').Value = 13000
$ws.Range('M215').Value = 12538
$ws.Range('N215').Value = '$/bandeja 18 kilos'
$ws.Range('O215').Value = 'Provincia de Limarí'
$ws.Range('P215').Value = 697
$ws.Range('Q215').Value = 18
$ws.Range('R215').Value = 'Hortaliza'
